$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 8209
$ws1.Range("G3").Value = 79.2
$ws1.Range("F4").Value = 1923
$ws1.Range("F5").Value = 6525
$ws1.Range("F7").Value = 2084
$ws1.Range("F11").Value = 56
$ws1.Range("F16").Value = 8582
$ws1.Range("F26").Value = 42
$ws1.Range("F32").Value = 8
$ws1.Range("F33").Value = 2120
$ws1.Range("F34").Value = 851
$ws1.Range("F35").Value = 491
$ws1.Range("F39").Value = 198
$ws1.Range("F42").Value = 45
$ws1.Range("F43").Value = 96

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 8209
$ws4.Range("G6").Value = 79.2
$ws4.Range("F9").Value = 1923
$ws4.Range("F10").Value = 6525
$ws4.Range("F11").Value = 2084
$ws4.Range("F17").Value = 56
$ws4.Range("F20").Value = 8582
$ws4.Range("F28").Value = 42
$ws4.Range("F33").Value = 8
$ws4.Range("F34").Value = 2120
$ws4.Range("F35").Value = 851
$ws4.Range("F37").Value = 491
$ws4.Range("F40").Value = 199
$ws4.Range("F43").Value = 96
